$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header area updates
# ---------------------------------------------------------------------------
# Total "Valor Mora" (E11)
$ws.Range("E11").Value = 600085
# "Cant. Trabajadores" count (C13)
$ws.Range("C13").Value = 9

# ---------------------------------------------------------------------------
# Grow the detail table by one row: insert a fresh row before the old last
# (bottom-bordered) row 24, copying the formatting of row 23 so the new row
# gets the regular inner-row style and the old row 24 keeps its bottom
# border style, now shifted down to row 25.
# ---------------------------------------------------------------------------
$ws.Rows("24:24").Insert()
$ws.Range("B23:J23").Copy()
$ws.Range("B24:J24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Replace the detail table contents (rows 16-25) with the new data set
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "45520291",   "DIANA YANETH FIGUEROA ARROYO",     "2507", 56940, 781242),
    @("CC", "45514109",   "MARLYS ESTHER ROMERO DIAZ",        "2507", 56940, 1423500),
    @("CC", "45524500",   "ANA MARGARITA ANGULO SANCHEZ",     "2507", 98172, 0),
    @("CC", "45536925",   "IRINA DEL CARMEN RIPOLL SANCHEZ",  "2507", 56940, 1423500),
    @("CC", "45536925",   "IRINA DEL CARMEN RIPOLL SANCHEZ",  "2205", 29333, 1423500),
    @("CC", "1047449477", "LEONARDO MIRANDA FAJARDO",         "2507", 56940, 1423500),
    @("CC", "1143358370", "CARLOS ANDRES ZUÑIGA GONZALEZ",    "2507", 56940, 908526),
    @("CC", "45586786",   "MARISOL MARRUGO BOSSIO",           "2507", 56940, 1423500),
    @("CC", "1047494192", "ANDREA PAOLA SINCELEJO JUNIELES",  "2507", 56940, 1423500),
    @("CC", "1047417996", "ALISSON CABEZA RAMOS",             "2507", 74000, 1850000)
)

$r = 16
foreach ($row in $data) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r = $r + 1
}
